$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting existing data (row 2 -> row 3) down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new user record.
$ws.Cells.Item(2, 1).Value = "dangoosh"
$ws.Cells.Item(2, 2).Value = "umrbek.xudayorovich@gmail.com"
$ws.Cells.Item(2, 3).Value = 1234567
$ws.Cells.Item(2, 4).Value = "2024-11-08T10:02:47.978767Z"
